$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2829.8
$ws.Range("J17").Value = 2500
$ws.Range("L17").Value = 7500
$ws.Range("N17").Value = -7836
$ws.Range("H33").Value = 544.2857
$ws.Range("I33").Value = 496.5263
$ws.Range("K33").Value = 496.5263
$ws.Range("M33").Value = -267.5263
$ws.Range("H34").Value = 16278.909
$ws.Range("I34").Value = 4340.8887
$ws.Range("K34").Value = 4340.8887
$ws.Range("M34").Value = -4137.8887
$ws.Range("H36").Value = 16278.909
$ws.Range("I36").Value = 4340.8887
$ws.Range("K36").Value = 4340.8887
$ws.Range("M36").Value = -3625.8887
$ws.Range("H43").Value = 6262.727
$ws.Range("J43").Value = 6871.143
$ws.Range("L43").Value = 6871.143
$ws.Range("N43").Value = -7009.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 3800
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H43").Value = 39997.5
$ws.Range("I43").Value = 39995
$ws.Range("K43").Value = 39995
$ws.Range("M43").Value = -39682
$ws.Range("H61").Value = 2012.8334
$ws.Range("I61").Value = 1612.6
$ws.Range("J61").Value = 4014
$ws.Range("K61").Value = 1612.6
$ws.Range("L61").Value = 4014
$ws.Range("M61").Value = -1400.6
$ws.Range("N61").Value = -4438
$ws.Range("H136").Value = 2012.8334
$ws.Range("I136").Value = 1612.6
$ws.Range("J136").Value = 4014
$ws.Range("K136").Value = 4837.799999999999
$ws.Range("L136").Value = 12042
$ws.Range("M136").Value = -2287.799999999999
$ws.Range("N136").Value = -17142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 22000
$ws.Range("J15").Value = 22000
$ws.Range("L15").Value = 22000
$ws.Range("N15").Value = -22454
$ws.Range("H20").Value = 392.5
$ws.Range("I20").Value = 392.5
$ws.Range("K20").Value = 392.5
$ws.Range("M20").Value = -145.5
$ws.Range("H22").Value = 714.8182
$ws.Range("I22").Value = 736.3
$ws.Range("K22").Value = 736.3
$ws.Range("M22").Value = -563.3
$ws.Range("H94").Value = 65627.64999999999
$ws.Range("I94").Value = 79483.57000000001
$ws.Range("K94").Value = 79483.57000000001
$ws.Range("M94").Value = -79032.57000000001
$ws.Range("H107").Value = 134650
$ws.Range("I107").Value = 134650
$ws.Range("K107").Value = 134650
$ws.Range("M107").Value = -132730

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 253.63637
$ws.Range("I7").Value = 265.70587
$ws.Range("J7").Value = 212.6
$ws.Range("K7").Value = 265.70587
$ws.Range("L7").Value = 212.6
$ws.Range("M7").Value = -152.70587
$ws.Range("N7").Value = -438.6
$ws.Range("H22").Value = 814.4
$ws.Range("I22").Value = 861.5
$ws.Range("J22").Value = 626
$ws.Range("K22").Value = 861.5
$ws.Range("L22").Value = 626
$ws.Range("M22").Value = -511.5
$ws.Range("N22").Value = -1326
$ws.Range("H31").Value = 5275.5
$ws.Range("I31").Value = 3195
$ws.Range("J31").Value = 6431.3335
$ws.Range("K31").Value = 3195
$ws.Range("L31").Value = 6431.3335
$ws.Range("M31").Value = -2900
$ws.Range("N31").Value = -7021.3335
$ws.Range("H32").Value = 3454
$ws.Range("I32").Value = 1817.5
$ws.Range("K32").Value = 1817.5
$ws.Range("M32").Value = -1501.5
$ws.Range("H34").Value = 5275.5
$ws.Range("I34").Value = 3195
$ws.Range("J34").Value = 6431.3335
$ws.Range("K34").Value = 3195
$ws.Range("L34").Value = 6431.3335
$ws.Range("M34").Value = -2993
$ws.Range("N34").Value = -6835.3335
$ws.Range("H50").Value = 23562.25
$ws.Range("J50").Value = 26900
$ws.Range("L50").Value = 26900
$ws.Range("N50").Value = -28150
$ws.Range("H99").Value = 456084.72
$ws.Range("I99").Value = 1717.3334
$ws.Range("K99").Value = 1717.3334
$ws.Range("M99").Value = -219.3334
$ws.Range("H126").Value = 456084.72
$ws.Range("I126").Value = 1717.3334
$ws.Range("K126").Value = 5152.0002
$ws.Range("M126").Value = -2682.0002
$ws.Range("H141").Value = 229888.2
$ws.Range("J141").Value = 427776.4
$ws.Range("L141").Value = 427776.4
$ws.Range("N141").Value = -438136.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 777.5
$ws.Range("I26").Value = 370
$ws.Range("K26").Value = 1110
$ws.Range("M26").Value = -822
$ws.Range("H113").Value = 687.2
$ws.Range("I113").Value = 711.5
$ws.Range("J113").Value = 590
$ws.Range("K113").Value = 2134.5
$ws.Range("L113").Value = 1770
$ws.Range("M113").Value = 35.5
$ws.Range("N113").Value = -6110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2064
$ws.Range("J11").Value = 128
$ws.Range("L11").Value = 128
$ws.Range("N11").Value = -406
$ws.Range("H18").Value = 802.5
$ws.Range("I18").Value = 802.5
$ws.Range("K18").Value = 802.5
$ws.Range("M18").Value = -509.5
$ws.Range("H22").Value = 3151.5
$ws.Range("J22").Value = 7332.6665
$ws.Range("L22").Value = 7332.6665
$ws.Range("N22").Value = -8390.666499999999
$ws.Range("H80").Value = 4500
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("H83").Value = 4500
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008
$ws.Range("H97").Value = 1122.8462
$ws.Range("I97").Value = 1122.8462
$ws.Range("K97").Value = 1122.8462
$ws.Range("M97").Value = -626.8462
$ws.Range("H105").Value = 9550
$ws.Range("J105").Value = 9550
$ws.Range("L105").Value = 9550
$ws.Range("N105").Value = -16538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 423
$ws.Range("I22").Value = 413
$ws.Range("J22").Value = 433
$ws.Range("K22").Value = 413
$ws.Range("L22").Value = 433
$ws.Range("M22").Value = -118
$ws.Range("N22").Value = -1023
$ws.Range("H27").Value = 423
$ws.Range("I27").Value = 413
$ws.Range("J27").Value = 433
$ws.Range("K27").Value = 413
$ws.Range("L27").Value = 433
$ws.Range("M27").Value = -306
$ws.Range("N27").Value = -647
$ws.Range("H55").Value = 977.2308
$ws.Range("J55").Value = 996.8570999999999
$ws.Range("L55").Value = 996.8570999999999
$ws.Range("N55").Value = -1342.8571
$ws.Range("H110").Value = 31999.5
$ws.Range("J110").Value = 31999.5
$ws.Range("L110").Value = 31999.5
$ws.Range("N110").Value = -40179.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 20000
$ws.Range("K53").Value = 20000
$ws.Range("M53").Value = -19393
$ws.Range("H62").Value = 2155.25
$ws.Range("I62").Value = 1999.5
$ws.Range("J62").Value = 2311
$ws.Range("K62").Value = 1999.5
$ws.Range("L62").Value = 2311
$ws.Range("M62").Value = -1375.5
$ws.Range("N62").Value = -3559
$ws.Range("H65").Value = 2155.25
$ws.Range("I65").Value = 1999.5
$ws.Range("J65").Value = 2311
$ws.Range("K65").Value = 9997.5
$ws.Range("L65").Value = 11555
$ws.Range("M65").Value = -6877.5
$ws.Range("N65").Value = -17795
$ws.Range("H113").Value = 364.5
$ws.Range("I113").Value = 202.3
$ws.Range("J113").Value = 634.8333
$ws.Range("K113").Value = 606.9000000000001
$ws.Range("L113").Value = 1904.4999
$ws.Range("M113").Value = 1563.1
$ws.Range("N113").Value = -6244.4999
